$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 277.42856
$ws.Range("I19").Value = 287.30768
$ws.Range("J19").Value = 268.86667
$ws.Range("K19").Value = 287.30768
$ws.Range("L19").Value = 268.86667
$ws.Range("M19").Value = -112.30768
$ws.Range("N19").Value = -618.86667
$ws.Range("H42").Value = 386.8889
$ws.Range("I42").Value = 172.4
$ws.Range("J42").Value = 655
$ws.Range("K42").Value = 517.2
$ws.Range("L42").Value = 1965
$ws.Range("M42").Value = -287.2
$ws.Range("N42").Value = -2425
$ws.Range("H129").Value = 187862.22
$ws.Range("I129").Value = 470.75
$ws.Range("J129").Value = 220452.05
$ws.Range("K129").Value = 1412.25
$ws.Range("L129").Value = 661356.1499999999
$ws.Range("M129").Value = 3587.75
$ws.Range("N129").Value = -671356.1499999999
$ws.Range("H138").Value = 3105.64
$ws.Range("I138").Value = 1473.898
$ws.Range("J138").Value = 4673.392
$ws.Range("K138").Value = 4421.694
$ws.Range("L138").Value = 14020.176
$ws.Range("M138").Value = 718.3060000000005
$ws.Range("N138").Value = -24300.176

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3518.18
$ws.Range("I32").Value = 3402.202
$ws.Range("K32").Value = 3402.202
$ws.Range("M32").Value = -3115.202
$ws.Range("H102").Value = 2299.2144
$ws.Range("I102").Value = 2280.7036
$ws.Range("K102").Value = 2280.7036
$ws.Range("M102").Value = -658.7035999999998
$ws.Range("H122").Value = 2317.2693
$ws.Range("I122").Value = 2059.2856
$ws.Range("K122").Value = 6177.8568
$ws.Range("M122").Value = -3727.8568
$ws.Range("H133").Value = 95697
$ws.Range("J133").Value = 95697
$ws.Range("L133").Value = 95697
$ws.Range("N133").Value = -100757

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1920
$ws.Range("I22").Value = 3147.5
$ws.Range("K22").Value = 3147.5
$ws.Range("M22").Value = -2797.5
$ws.Range("H122").Value = 1671417.4
$ws.Range("J122").Value = 2005098.4
$ws.Range("L122").Value = 6015295.199999999
$ws.Range("N122").Value = -6020195.199999999
$ws.Range("H134").Value = 1995.8158
$ws.Range("I134").Value = 2274.8
$ws.Range("J134").Value = 1459.3077
$ws.Range("K134").Value = 6824.400000000001
$ws.Range("L134").Value = 4377.9231
$ws.Range("M134").Value = -4289.400000000001
$ws.Range("N134").Value = -9447.9231

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1229.875
$ws.Range("J113").Value = 509.1111
$ws.Range("L113").Value = 1527.3333
$ws.Range("N113").Value = -5867.3333
$ws.Range("H127").Value = 1031.6666
$ws.Range("J127").Value = 1031.6666
$ws.Range("L127").Value = 3094.9998
$ws.Range("N127").Value = -13014.9998
$ws.Range("H131").Value = 847.92
$ws.Range("J131").Value = 870.1383
$ws.Range("L131").Value = 2610.4149
$ws.Range("N131").Value = -12690.4149

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1265.8334
$ws.Range("I113").Value = 1152.8823
$ws.Range("J113").Value = 1540.1428
$ws.Range("K113").Value = 1152.8823
$ws.Range("L113").Value = 1540.1428
$ws.Range("M113").Value = 1017.1177
$ws.Range("N113").Value = -5880.1428
$ws.Range("H122").Value = 2499.325
$ws.Range("I122").Value = 2043.25
$ws.Range("J122").Value = 3563.5
$ws.Range("K122").Value = 6129.75
$ws.Range("L122").Value = 10690.5
$ws.Range("M122").Value = -3679.75
$ws.Range("N122").Value = -15590.5
$ws.Range("H132").Value = 1961.7556
$ws.Range("I132").Value = 1571.0646
$ws.Range("J132").Value = 2826.8572
$ws.Range("K132").Value = 4713.1938
$ws.Range("L132").Value = 8480.571599999999
$ws.Range("M132").Value = -2183.1938
$ws.Range("N132").Value = -13540.5716
$ws.Range("H136").Value = 9597.183999999999
$ws.Range("J136").Value = 9597.183999999999
$ws.Range("L136").Value = 28791.552
$ws.Range("N136").Value = -33891.552

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 62677.883
$ws.Range("I7").Value = 94429.55
$ws.Range("J7").Value = 4466.5
$ws.Range("K7").Value = 94429.55
$ws.Range("L7").Value = 4466.5
$ws.Range("M7").Value = -94317.55
$ws.Range("N7").Value = -4690.5
$ws.Range("H22").Value = 1143.5454
$ws.Range("I22").Value = 694.75
$ws.Range("J22").Value = 1400
$ws.Range("K22").Value = 694.75
$ws.Range("L22").Value = 1400
$ws.Range("M22").Value = -399.75
$ws.Range("N22").Value = -1990
$ws.Range("H27").Value = 1143.5454
$ws.Range("I27").Value = 694.75
$ws.Range("J27").Value = 1400
$ws.Range("K27").Value = 694.75
$ws.Range("L27").Value = 1400
$ws.Range("M27").Value = -587.75
$ws.Range("N27").Value = -1614
$ws.Range("H93").Value = 776.9
$ws.Range("I93").Value = 776.9
$ws.Range("K93").Value = 776.9
$ws.Range("M93").Value = 471.1
$ws.Range("H122").Value = 7938608.5
$ws.Range("I122").Value = 13889992
$ws.Range("J122").Value = 3430
$ws.Range("K122").Value = 41669976
$ws.Range("L122").Value = 10290
$ws.Range("M122").Value = -41667526
$ws.Range("N122").Value = -15190
$ws.Range("H126").Value = 62677.883
$ws.Range("I126").Value = 94429.55
$ws.Range("J126").Value = 4466.5
$ws.Range("K126").Value = 283288.65
$ws.Range("L126").Value = 13399.5
$ws.Range("M126").Value = -280818.65
$ws.Range("N126").Value = -18339.5
$ws.Range("H132").Value = 7379.3286
$ws.Range("I132").Value = 6127.7593
$ws.Range("J132").Value = 10936.421
$ws.Range("K132").Value = 18383.2779
$ws.Range("L132").Value = 32809.263
$ws.Range("M132").Value = -15853.2779
$ws.Range("N132").Value = -37869.263

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1846.8966
$ws.Range("I81").Value = 897.1
$ws.Range("J81").Value = 2346.7896
$ws.Range("K81").Value = 1794.2
$ws.Range("L81").Value = 4693.5792
$ws.Range("M81").Value = -733.2
$ws.Range("N81").Value = -6815.5792
$ws.Range("H84").Value = 1846.8966
$ws.Range("I84").Value = 897.1
$ws.Range("J84").Value = 2346.7896
$ws.Range("K84").Value = 8971
$ws.Range("L84").Value = 23467.896
$ws.Range("M84").Value = -3667
$ws.Range("N84").Value = -34075.896
$ws.Range("H96").Value = 1700
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 1700
$ws.Range("K96").Value = 0
$ws.Range("L96").ClearContents()
$ws.Range("M96").Value = 1700
$ws.Range("N96").Value = -4446
$ws.Range("H113").Value = 200.18182
$ws.Range("I113").Value = 197.1875
$ws.Range("K113").Value = 591.5625
$ws.Range("M113").Value = 1578.4375
$ws.Range("H122").Value = 49659.285
$ws.Range("I122").Value = 64340
$ws.Range("J122").Value = 2681
$ws.Range("K122").Value = 193020
$ws.Range("L122").Value = 8043
$ws.Range("M122").Value = -190570
$ws.Range("N122").Value = -12943
$ws.Range("H132").Value = 1206.035
$ws.Range("I132").Value = 900.6977000000001
$ws.Range("J132").Value = 2143.8572
$ws.Range("K132").Value = 2702.0931
$ws.Range("L132").Value = 6431.571599999999
$ws.Range("M132").Value = -172.0931
$ws.Range("N132").Value = -11491.5716
